# Update cryptos list with latest price/volume data (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 47/48: NEARProtocol and PancakeSwap swapped ranking positions.
# Leading "'" keeps numeric-looking Price values stored as text (matching
# the source data's inline-string cells) instead of being auto-converted
# to numbers by Excel, which would drop significant trailing zeros.
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.906"
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.233"
$ws.Range("E48").Value = "  -3.73%  "

# Price (D) / Volume(1h) (E) updates for all other rows.
$ws.Range("D2").Value = "23.383.11"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "1.632.79"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "'299.97"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("D7").Value = "'0.3777"
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("D9").Value = "'0.3521"
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("D10").Value = "'0.08042"
$ws.Range("E10").Value = "  -2.08%  "
$ws.Range("E11").Value = "  -3.60%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D14").Value = "'6.298"
$ws.Range("E14").Value = "  -3.34%  "
$ws.Range("D15").Value = "'7.230"
$ws.Range("E15").Value = "  -2.81%  "
$ws.Range("D16").Value = "'0.00001199"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").Value = "1.630.05"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Value = "'95.58"
$ws.Range("E18").Value = "  -1.84%  "
$ws.Range("D19").Value = "'0.06953"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "'6.667"
$ws.Range("E20").Value = "  -2.01%  "
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'12.28"
$ws.Range("D24").Value = "23.391.47"
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("D25").Value = "'2.459"
$ws.Range("E25").Value = "  -2.72%  "
$ws.Range("D26").Value = "'2.880"
$ws.Range("E26").Value = "  -5.65%  "
$ws.Range("D27").Value = "'20.75"
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("D28").Value = "'152.10"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "'5.180"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").Value = "'132.06"
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("D31").Value = "1.813.24"
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("D32").Value = "'6.797"
$ws.Range("E32").Value = "  -1.93%  "
$ws.Range("D33").Value = "'2.129"
$ws.Range("E33").Value = "  -3.30%  "
$ws.Range("D34").Value = "'11.29"
$ws.Range("E34").Value = "  -4.08%  "
$ws.Range("D35").Value = "'0.9676"
$ws.Range("E35").Value = "  -9.79%  "
$ws.Range("D36").Value = "'0.02704"
$ws.Range("E36").Value = "  -3.84%  "
$ws.Range("D37").Value = "'0.08714"
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("D38").Value = "'0.2429"
$ws.Range("E38").Value = "  -4.15%  "
$ws.Range("D39").Value = "'5.858"
$ws.Range("E39").Value = "  -4.00%  "
$ws.Range("D40").Value = "'0.06795"
$ws.Range("E40").Value = "  -4.63%  "
$ws.Range("D41").Value = "'12.85"
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("D42").Value = "'0.6806"
$ws.Range("E42").Value = "  -3.16%  "
$ws.Range("D43").Value = "'1.295"
$ws.Range("E43").Value = "  -3.26%  "
$ws.Range("D44").Value = "'15.39"
$ws.Range("E44").Value = "  -4.21%  "
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "'0.6297"
$ws.Range("E46").Value = "  -3.23%  "
$ws.Range("D49").Value = "'0.07670"
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("D50").Value = "'126.53"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").Value = "'1.199"
$ws.Range("E51").Value = "  +0.32%  "
